$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (e.g. "299.10" -> 299.1).
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D17",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values as described by the diff.
$ws.Range("D2").Value = '44.327.46'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '2.224.92'
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '299.10'
$ws.Range("E5").Value = '  -2.99%  '
$ws.Range("D6").Value = '89.88'
$ws.Range("E6").Value = '  -5.57%  '
$ws.Range("D7").Value = '0.561'
$ws.Range("E7").Value = '  -2.16%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  -6.85%  '
$ws.Range("D10").Value = '32.93'
$ws.Range("E10").Value = '  -6.82%  '
$ws.Range("D11").Value = '0.0782'
$ws.Range("E11").Value = '  -3.44%  '
$ws.Range("D12").Value = '6.95'
$ws.Range("E12").Value = '  -4.95%  '
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").Value = '2.565.02'
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").Value = '2.212.06'
$ws.Range("E15").Value = '  -4.17%  '
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("D17").Value = '0.779'
$ws.Range("E17").Value = '  -7.40%  '
$ws.Range("D18").Value = '44.104.79'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = '0.0₃0911'
$ws.Range("E19").Value = '  -5.81%  '
$ws.Range("D20").Value = '5.91'
$ws.Range("E20").Value = '  -8.08%  '
$ws.Range("D21").Value = '10.99'
$ws.Range("E21").Value = '  -10.26%  '
$ws.Range("D22").Value = '64.66'
$ws.Range("E22").Value = '  -2.12%  '
$ws.Range("D23").Value = '238.25'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '2.80'
$ws.Range("E24").Value = '  -6.63%  '
$ws.Range("E25").Value = '  +0.94%  '
$ws.Range("D26").Value = '1.86'
$ws.Range("E26").Value = '  -7.64%  '
$ws.Range("D27").Value = '2.23'
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").Value = '38.37'
$ws.Range("E28").Value = '  +0.82%  '
$ws.Range("D29").Value = '9.35'
$ws.Range("E29").Value = '  -5.56%  '
$ws.Range("D30").Value = '19.49'
$ws.Range("E30").Value = '  -3.12%  '
$ws.Range("D31").Value = '148.41'
$ws.Range("E31").Value = '  -2.79%  '
$ws.Range("D32").Value = '5.39'
$ws.Range("E32").Value = '  -10.28%  '
$ws.Range("E33").Value = '  -3.58%  '
$ws.Range("D34").Value = '0.0749'
$ws.Range("E34").Value = '  -6.87%  '
$ws.Range("E35").Value = '  -4.10%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.103'
$ws.Range("E36").Value = '  -5.98%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '2.81'
$ws.Range("E37").Value = '  -11.77%  '
$ws.Range("D38").Value = '1.68'
$ws.Range("E38").Value = '  -7.13%  '
$ws.Range("D39").Value = '0.0302'
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("D40").Value = '3.19'
$ws.Range("E40").Value = '  -7.64%  '
$ws.Range("D41").Value = '3.52'
$ws.Range("E41").Value = '  -8.11%  '
$ws.Range("D42").Value = '13.13'
$ws.Range("E42").Value = '  -10.11%  '
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("D44").Value = '1.810.88'
$ws.Range("E44").Value = '  +3.34%  '
$ws.Range("D45").Value = '1.77'
$ws.Range("E45").Value = '  +11.21%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.178'
$ws.Range("E46").Value = '  -8.28%  '
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").Value = '74.27'
$ws.Range("E47").Value = '  -8.45%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '14.10'
$ws.Range("E48").Value = '  +8.74%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '94.04'
$ws.Range("E49").Value = '  -6.34%  '
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").Value = '66.70'
$ws.Range("E50").Value = '  -6.31%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.446.64'
$ws.Range("E51").Value = '  -1.08%  '
